# Generate Report for Archive
# Re-sorts the "9f897edb..." / "9170348f..." rows (row 6 and row 7) into
# alphabetical order on every sheet (Overview, zh-cn, de-de), matching the
# canonical localization-status report generator output.

$wb = $excel.ActiveWorkbook

function Update-Hyperlink($ws, $row, $col, $newDisplay) {
    foreach ($link in $ws.Hyperlinks) {
        if (($link.Range.Row -eq $row) -and ($link.Range.Column -eq $col)) {
            $link.TextToDisplay = $newDisplay
        }
    }
}

# ---- Overview sheet: columns A (file name), B (zh-cn status), C (de-de status) ----
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A6").Value = "9170348f-dbc4-4623-bda6-b768f2ef8475.md"
$ws.Range("A7").Value = "9f897edb-7a44-463d-be91-d41781ed7fbc.md"
$ws.Range("B7").Value = "In Translation"
$ws.Range("C7").Value = "In Translation"

Update-Hyperlink $ws 6 1 "9170348f-dbc4-4623-bda6-b768f2ef8475.md"
Update-Hyperlink $ws 7 1 "9f897edb-7a44-463d-be91-d41781ed7fbc.md"

# ---- Per-language detail sheets: zh-cn and de-de ----
$langSheets = @(
    @{ Name = "zh-cn";
       C6 = "9f897edb-7a44-463d-be91-d41781ed7fbc.eb5938e2d7b4cb2a8c31cb133bb98818cfe68680.zh-cn.xlf";
       D6 = "2016-03-08 10:25:38";
       C7 = "9170348f-dbc4-4623-bda6-b768f2ef8475.edcab04bba996e067a9ce5d4280fcb16bd419568.zh-cn.xlf";
       D7 = "2016-03-08 10:28:01" },
    @{ Name = "de-de";
       C6 = "9f897edb-7a44-463d-be91-d41781ed7fbc.eb5938e2d7b4cb2a8c31cb133bb98818cfe68680.de-de.xlf";
       D6 = "2016-03-08 10:25:42";
       C7 = "9170348f-dbc4-4623-bda6-b768f2ef8475.edcab04bba996e067a9ce5d4280fcb16bd419568.de-de.xlf";
       D7 = "2016-03-08 10:28:05" }
)

foreach ($entry in $langSheets) {
    $ws = $wb.Worksheets.Item($entry.Name)

    # New row 6 becomes the former row 7 ("9170348f...") content
    $ws.Range("A6").Value = "9170348f-dbc4-4623-bda6-b768f2ef8475.md"
    $ws.Range("B6").Value = "In Translation"
    $ws.Range("C6").Value = $entry.C7
    $ws.Range("D6").Value = $entry.D7

    # New row 7 becomes the former row 6 ("9f897edb...") content
    $ws.Range("A7").Value = "9f897edb-7a44-463d-be91-d41781ed7fbc.md"
    $ws.Range("B7").Value = "In Translation"
    $ws.Range("C7").Value = $entry.C6
    $ws.Range("D7").Value = $entry.D6

    Update-Hyperlink $ws 6 1 "9170348f-dbc4-4623-bda6-b768f2ef8475.md"
    Update-Hyperlink $ws 6 3 $entry.C7
    Update-Hyperlink $ws 7 1 "9f897edb-7a44-463d-be91-d41781ed7fbc.md"
    Update-Hyperlink $ws 7 3 $entry.C6
}
